$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header, matching the bold style already used by A1/B1
$ws.Range("C1").Value = "League Position 2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# League position values for each team row (rows 2-8), in the same
# Westham/Arsenal/Chelsea/Man Utd/Newcastle/Wolves/Leeds order as
# columns A and B
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 12
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 4
$ws.Range("C7").Value = 13
$ws.Range("C8").Value = 19

# Average formula row
$ws.Range("C9").Formula = "=AVERAGE(C2:C8)"

$wb.Save()
